$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This edit inserts a new feature-description row ("nmhc_gt") after row 3,
# and appends two new rows ("month", "hour") after the existing table, while
# updating the explanatory text in column H to keep it aligned with column A.
#
# To avoid perturbing cells that are not supposed to change (and to avoid
# floating point re-serialization noise on untouched cells), we do NOT use
# Rows.Insert() to shift the grid. Instead we:
#   1. Leave row 1 (headers) and the unaffected A:G values of rows 2-3 alone.
#   2. Directly (re)write every cell of rows 4-16 with its final target
#      value (new row plus the old rows 4-13 shifted down by one, plus the
#      two brand new trailing rows), only touching H2/H3 for the two rows
#      whose underlying data doesn't move.
#   3. Copy formatting (font/border/alignment) for the new column-A cells
#      from an existing styled cell (A2) *before* writing the new value, so
#      the appended rows look the same as the rest of the table (bold,
#      bordered, centered) without leaving unused style entries behind.
#
# The order in which *new* strings are first written controls the order in
# which they are appended to the shared string table, so we deliberately
# pre-touch the brand new strings (via a scratch cell) in the same relative
# order they appear in the final workbook: "nmhc_gt", "month", "hour", then
# the long German explanation strings.
# ---------------------------------------------------------------------------

# Pre-seed the shared string table so brand-new strings are appended in the
# same relative order they appear in the final workbook (the engine appends
# newly-seen unique strings to the shared string table in first-touch
# order). We use a scratch cell far outside the used range and clear it
# afterwards.
$scratch = $ws.Cells.Item(1, 100)
$scratch.Value2 = "nmhc_gt"
$scratch.Value2 = "month"
$scratch.Value2 = "hour"
$scratch.Value2 = "Stuendlich gemittelte Gesamtkonzentration an nicht-metanischem Kohlenwasserstoff"
$scratch.Value2 = "Monate der Erfassung"
$scratch.Value2 = "Stunden der erfassung"
$scratch.ClearContents()

# --- Update H2 / H3 (explanation column) to the shifted shared strings.
#     Column A and B:G of rows 2-3 are intentionally left untouched. ---
$ws.Cells.Item(2, 8).Value2 = "Stuendlich gemittelte CO-Konzentration"
$ws.Cells.Item(3, 8).Value2 = "Stuendlich gemittelte Sensorreaktion (nominell auf CO ausgerichtet) (Zinnoxid)"

# --- New row 4: nmhc_gt (inserted feature) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(4, 1))
$ws.Cells.Item(4, 1).Value2 = "nmhc_gt"
$ws.Cells.Item(4, 2).Value2 = 218.6076662908681
$ws.Cells.Item(4, 3).Value2 = 7
$ws.Cells.Item(4, 4).Value2 = 66
$ws.Cells.Item(4, 5).Value2 = 145
$ws.Cells.Item(4, 6).Value2 = 297
$ws.Cells.Item(4, 7).Value2 = 1189
$ws.Cells.Item(4, 8).Value2 = "Stuendlich gemittelte Gesamtkonzentration an nicht-metanischem Kohlenwasserstoff"

# --- Row 5: c6h6_gt (was row 4) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(5, 1))
$ws.Cells.Item(5, 1).Value2 = "c6h6_gt"
$ws.Cells.Item(5, 2).Value2 = 10.4582046958101
$ws.Cells.Item(5, 3).Value2 = 0.1
$ws.Cells.Item(5, 4).Value2 = 4.7
$ws.Cells.Item(5, 5).Value2 = 8.6
$ws.Cells.Item(5, 6).Value2 = 14.4
$ws.Cells.Item(5, 7).Value2 = 63.7
$ws.Cells.Item(5, 8).Value2 = "Stuendlich gemittelte Benzolkonzentration"

# --- Row 6: pt08_s2_nmhc (was row 5) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(6, 1))
$ws.Cells.Item(6, 1).Value2 = "pt08_s2_nmhc"
$ws.Cells.Item(6, 2).Value2 = 953.5794525878842
$ws.Cells.Item(6, 3).Value2 = 383
$ws.Cells.Item(6, 4).Value2 = 749
$ws.Cells.Item(6, 5).Value2 = 925
$ws.Cells.Item(6, 6).Value2 = 1130
$ws.Cells.Item(6, 7).Value2 = 2214
$ws.Cells.Item(6, 8).Value2 = "Stuendlich gemittelte Sensorreaktion (nominell auf NMHC ausgerichtet) (Titandioxid)"

# --- Row 7: nox_gt (was row 6) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(7, 1))
$ws.Cells.Item(7, 1).Value2 = "nox_gt"
$ws.Cells.Item(7, 2).Value2 = 232.3600715912789
$ws.Cells.Item(7, 3).Value2 = 2
$ws.Cells.Item(7, 4).Value2 = 89
$ws.Cells.Item(7, 5).Value2 = 164
$ws.Cells.Item(7, 6).Value2 = 303.75
$ws.Cells.Item(7, 7).Value2 = 1479
$ws.Cells.Item(7, 8).Value2 = "Echte Stuendlich gemittelte NOx-Konzentration"

# --- Row 8: pt08_s3_nox (was row 7) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(8, 1))
$ws.Cells.Item(8, 1).Value2 = "pt08_s3_nox"
$ws.Cells.Item(8, 2).Value2 = 848.7626151251784
$ws.Cells.Item(8, 3).Value2 = 322
$ws.Cells.Item(8, 4).Value2 = 672
$ws.Cells.Item(8, 5).Value2 = 818
$ws.Cells.Item(8, 6).Value2 = 984
$ws.Cells.Item(8, 7).Value2 = 2683
$ws.Cells.Item(8, 8).Value2 = "Suendlich gemitteltes Sensoransprechverhalten (nominell auf NOx ausgerichtet)"

# --- Row 9: no2_gt (was row 8) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(9, 1))
$ws.Cells.Item(9, 1).Value2 = "no2_gt"
$ws.Cells.Item(9, 2).Value2 = 106.2518313527592
$ws.Cells.Item(9, 3).Value2 = 2
$ws.Cells.Item(9, 4).Value2 = 73
$ws.Cells.Item(9, 5).Value2 = 103
$ws.Cells.Item(9, 6).Value2 = 132
$ws.Cells.Item(9, 7).Value2 = 333
$ws.Cells.Item(9, 8).Value2 = "Stuendlich gemittelte NO2-Konzentration"

# --- Row 10: pt08_s4_no2 (was row 9) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(10, 1))
$ws.Cells.Item(10, 1).Value2 = "pt08_s4_no2"
$ws.Cells.Item(10, 2).Value2 = 1507.352834349462
$ws.Cells.Item(10, 3).Value2 = 657
$ws.Cells.Item(10, 4).Value2 = 1305
$ws.Cells.Item(10, 5).Value2 = 1508
$ws.Cells.Item(10, 6).Value2 = 1707
$ws.Cells.Item(10, 7).Value2 = 2775
$ws.Cells.Item(10, 8).Value2 = "Stuendlich gemittelte Sensorreaktion (nominell auf NO2 ausgerichtet) (Wolframoxid)"

# --- Row 11: pt08_s5_o3 (was row 10) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(11, 1))
$ws.Cells.Item(11, 1).Value2 = "pt08_s5_o3"
$ws.Cells.Item(11, 2).Value2 = 1024.200025943702
$ws.Cells.Item(11, 3).Value2 = 253
$ws.Cells.Item(11, 4).Value2 = 737
$ws.Cells.Item(11, 5).Value2 = 962
$ws.Cells.Item(11, 6).Value2 = 1272
$ws.Cells.Item(11, 7).Value2 = 2523
$ws.Cells.Item(11, 8).Value2 = "Stuendlich gemitteltes Sensoransprechverhalten (nominell O3-bezogen) (Indiumoxid)"

# --- Row 12: t (was row 11) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(12, 1))
$ws.Cells.Item(12, 1).Value2 = "t"
$ws.Cells.Item(12, 2).Value2 = 19.47604099104944
$ws.Cells.Item(12, 3).Value2 = 0.3
$ws.Cells.Item(12, 4).Value2 = 13.1
$ws.Cells.Item(12, 5).Value2 = 19.3
$ws.Cells.Item(12, 6).Value2 = 25.4
$ws.Cells.Item(12, 7).Value2 = 44.6
$ws.Cells.Item(12, 8).Value2 = "Temperatur"

# --- Row 13: rh (was row 12) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(13, 1))
$ws.Cells.Item(13, 1).Value2 = "rh"
$ws.Cells.Item(13, 2).Value2 = 48.96924374108173
$ws.Cells.Item(13, 3).Value2 = 9.199999999999999
$ws.Cells.Item(13, 4).Value2 = 35.5
$ws.Cells.Item(13, 5).Value2 = 49.4
$ws.Cells.Item(13, 6).Value2 = 62.1
$ws.Cells.Item(13, 7).Value2 = 88.7
$ws.Cells.Item(13, 8).Value2 = "Relative Luftfeuchtigkeit"

# --- Row 14: ah (was row 13) ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(14, 1))
$ws.Cells.Item(14, 1).Value2 = "ah"
$ws.Cells.Item(14, 2).Value2 = 1.082098624983781
$ws.Cells.Item(14, 3).Value2 = 0.1988
$ws.Cells.Item(14, 4).Value2 = 0.8167
$ws.Cells.Item(14, 5).Value2 = 1.0468
$ws.Cells.Item(14, 6).Value2 = 1.3713
$ws.Cells.Item(14, 7).Value2 = 2.231
$ws.Cells.Item(14, 8).Value2 = "Absolute Luftfeuchtigkeit"

# --- New row 15: month ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(15, 1))
$ws.Cells.Item(15, 1).Value2 = "month"
$ws.Cells.Item(15, 2).Value2 = 6.928525100531846
$ws.Cells.Item(15, 3).Value2 = 1
$ws.Cells.Item(15, 4).Value2 = 4
$ws.Cells.Item(15, 5).Value2 = 7
$ws.Cells.Item(15, 6).Value2 = 10
$ws.Cells.Item(15, 7).Value2 = 12
$ws.Cells.Item(15, 8).Value2 = "Monate der Erfassung"

# --- New row 16: hour ---
$ws.Cells.Item(2, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(16, 1).Value2 = "hour"
$ws.Cells.Item(16, 2).Value2 = 11.48851991179141
$ws.Cells.Item(16, 3).Value2 = 0
$ws.Cells.Item(16, 4).Value2 = 6
$ws.Cells.Item(16, 5).Value2 = 11
$ws.Cells.Item(16, 6).Value2 = 17
$ws.Cells.Item(16, 7).Value2 = 23
$ws.Cells.Item(16, 8).Value2 = "Stunden der erfassung"
